$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.345.33"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "3.538.62"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.12"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.27"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D7").Value = "3.536.07"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "8.15"
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("E11").Value = "  -4.12%  "
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("D13").Value = "4.138.35"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E14").Value = "  -4.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.30"
$ws.Range("E15").Value = "  -5.31%  "
$ws.Range("D16").Value = "3.536.50"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "66.407.42"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.99"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.94"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.23"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.602"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.86"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Value = "3.676.58"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -5.78%  "
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -7.13%  "
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.31"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("D35").Value = "3.528.54"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.83"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.64"
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.32"
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0857"
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.23"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.894"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").Value = "  -7.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.74"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.22"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.13"
$ws.Range("E48").Value = "  -7.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.41"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  -4.69%  "
$ws.Range("E51").Value = "  -5.38%  "
